$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Cells.Item(9, 8).Value2 = 135.93333  # H9: 133.6875 -> 135.93333
$ws.Cells.Item(9, 9).Value2 = 111.1  # I9: 110.09091 -> 111.1
$ws.Cells.Item(9, 11).Value2 = 111.1  # K9: 110.09091 -> 111.1
$ws.Cells.Item(9, 13).Value2 = 57.90000000000001  # M9: 58.90909000000001 -> 57.90000000000001
# Row 33
$ws.Cells.Item(33, 8).Value2 = 463.875  # H33: 489.26666 -> 463.875
$ws.Cells.Item(33, 9).Value2 = 454.8  # I33: 481.35715 -> 454.8
$ws.Cells.Item(33, 11).Value2 = 454.8  # K33: 481.35715 -> 454.8
$ws.Cells.Item(33, 13).Value2 = -225.8  # M33: -252.35715 -> -225.8
# Row 40
$ws.Cells.Item(40, 8).Value2 = 1001  # H40: 1667.3334 -> 1001
$ws.Cells.Item(40, 10).Value2 = 0  # J40: 3000 -> 0
$ws.Cells.Item(40, 12).Value2 = 0  # L40: 3000 -> 0
$ws.Cells.Item(40, 14).ClearContents()  # N40: was -3350
# Row 64
$ws.Cells.Item(64, 8).Value2 = 4150  # H64: 3750.6 -> 4150
$ws.Cells.Item(64, 10).Value2 = 3730  # J64: 3440.75 -> 3730
$ws.Cells.Item(64, 12).Value2 = 3730  # L64: 3440.75 -> 3730
$ws.Cells.Item(64, 14).Value2 = -4226  # N64: -3936.75 -> -4226
# Row 67
$ws.Cells.Item(67, 8).Value2 = 4150  # H67: 3750.6 -> 4150
$ws.Cells.Item(67, 10).Value2 = 3730  # J67: 3440.75 -> 3730
$ws.Cells.Item(67, 12).Value2 = 3730  # L67: 3440.75 -> 3730
$ws.Cells.Item(67, 14).Value2 = -5446  # N67: -5156.75 -> -5446
# Row 76
$ws.Cells.Item(76, 8).Value2 = 5012  # H76: 5566 -> 5012
$ws.Cells.Item(76, 9).Value2 = 3774.5  # I76: 4199 -> 3774.5
$ws.Cells.Item(76, 11).Value2 = 3774.5  # K76: 4199 -> 3774.5
$ws.Cells.Item(76, 13).Value2 = -3459.5  # M76: -3884 -> -3459.5
# Row 79
$ws.Cells.Item(79, 8).Value2 = 5012  # H79: 5566 -> 5012
$ws.Cells.Item(79, 9).Value2 = 3774.5  # I79: 4199 -> 3774.5
$ws.Cells.Item(79, 11).Value2 = 3774.5  # K79: 4199 -> 3774.5
$ws.Cells.Item(79, 13).Value2 = -2682.5  # M79: -3107 -> -2682.5
# Row 98
$ws.Cells.Item(98, 8).Value2 = 4012.4062  # H98: 4012.7188 -> 4012.4062
$ws.Cells.Item(98, 9).Value2 = 4266.9  # I98: 4267.2334 -> 4266.9
$ws.Cells.Item(98, 11).Value2 = 4266.9  # K98: 4267.2334 -> 4266.9
$ws.Cells.Item(98, 13).Value2 = -2768.9  # M98: -2769.2334 -> -2768.9
# Row 116
$ws.Cells.Item(116, 8).Value2 = 2943.2856  # H116: 2603.4285 -> 2943.2856
$ws.Cells.Item(116, 9).Value2 = 2868  # I116: 2356.125 -> 2868
$ws.Cells.Item(116, 10).Value2 = 2999.75  # J116: 2933.1667 -> 2999.75
$ws.Cells.Item(116, 11).Value2 = 2868  # K116: 2356.125 -> 2868
$ws.Cells.Item(116, 12).Value2 = 2999.75  # L116: 2933.1667 -> 2999.75
$ws.Cells.Item(116, 13).Value2 = 574  # M116: 1085.875 -> 574
$ws.Cells.Item(116, 14).Value2 = -9883.75  # N116: -9817.1667 -> -9883.75
# Row 122
$ws.Cells.Item(122, 8).Value2 = 4012.4062  # H122: 4012.7188 -> 4012.4062
$ws.Cells.Item(122, 9).Value2 = 4266.9  # I122: 4267.2334 -> 4266.9
$ws.Cells.Item(122, 11).Value2 = 12800.7  # K122: 12801.7002 -> 12800.7
$ws.Cells.Item(122, 13).Value2 = -10350.7  # M122: -10351.7002 -> -10350.7
# Row 132
$ws.Cells.Item(132, 8).Value2 = 5958332  # H132: 5853856 -> 5958332
$ws.Cells.Item(132, 9).Value2 = 6539680  # I132: 6670469.5 -> 6539680
$ws.Cells.Item(132, 10).Value2 = 28581.2  # J132: 20901.715 -> 28581.2
$ws.Cells.Item(132, 11).Value2 = 19619040  # K132: 20011408.5 -> 19619040
$ws.Cells.Item(132, 12).Value2 = 85743.60000000001  # L132: 62705.145 -> 85743.60000000001
$ws.Cells.Item(132, 13).Value2 = -19616510  # M132: -20008878.5 -> -19616510
$ws.Cells.Item(132, 14).Value2 = -90803.60000000001  # N132: -67765.145 -> -90803.60000000001
# Row 136
$ws.Cells.Item(136, 8).Value2 = 41747.273  # H136: 41861.816 -> 41747.273
$ws.Cells.Item(136, 10).Value2 = 41747.273  # J136: 41861.816 -> 41747.273
$ws.Cells.Item(136, 12).Value2 = 41747.273  # L136: 41861.816 -> 41747.273
$ws.Cells.Item(136, 14).Value2 = -51947.273  # N136: -52061.816 -> -51947.273
# Row 138
$ws.Cells.Item(138, 8).Value2 = 2307.34  # H138: 2185.31 -> 2307.34
$ws.Cells.Item(138, 9).Value2 = 1684.4375  # I138: 968.7857 -> 1684.4375
$ws.Cells.Item(138, 10).Value2 = 2425.988  # J138: 2383.3489 -> 2425.988
$ws.Cells.Item(138, 11).Value2 = 5053.3125  # K138: 2906.3571 -> 5053.3125
$ws.Cells.Item(138, 12).Value2 = 7277.964  # L138: 7150.0467 -> 7277.964
$ws.Cells.Item(138, 13).Value2 = 86.6875  # M138: 2233.6429 -> 86.6875
$ws.Cells.Item(138, 14).Value2 = -17557.964  # N138: -17430.0467 -> -17557.964

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value2 = 8732.382  # H32: 8988.893 -> 8732.382
$ws.Cells.Item(32, 9).Value2 = 6324.3296  # I32: 6499.926 -> 6324.3296
$ws.Cells.Item(32, 11).Value2 = 6324.3296  # K32: 6499.926 -> 6324.3296
$ws.Cells.Item(32, 13).Value2 = -6037.3296  # M32: -6212.926 -> -6037.3296
# Row 61
$ws.Cells.Item(61, 8).Value2 = 1593.7368  # H61: 66668280 -> 1593.7368
$ws.Cells.Item(61, 9).Value2 = 1493.7142  # I61: 76924510 -> 1493.7142
$ws.Cells.Item(61, 10).Value2 = 1873.8  # J61: 2749.5 -> 1873.8
$ws.Cells.Item(61, 11).Value2 = 1493.7142  # K61: 76924510 -> 1493.7142
$ws.Cells.Item(61, 12).Value2 = 1873.8  # L61: 2749.5 -> 1873.8
$ws.Cells.Item(61, 13).Value2 = -1281.7142  # M61: -76924298 -> -1281.7142
$ws.Cells.Item(61, 14).Value2 = -2297.8  # N61: -3173.5 -> -2297.8
# Row 97
$ws.Cells.Item(97, 8).Value2 = 526.5161000000001  # H97: 542.4 -> 526.5161000000001
$ws.Cells.Item(97, 9).Value2 = 501.85184  # I97: 519.2308 -> 501.85184
$ws.Cells.Item(97, 11).Value2 = 501.85184  # K97: 519.2308 -> 501.85184
$ws.Cells.Item(97, 13).Value2 = -5.851839999999982  # M97: -23.23080000000004 -> -5.851839999999982
# Row 132
$ws.Cells.Item(132, 8).Value2 = 2396.2727  # H132: 2616.2708 -> 2396.2727
$ws.Cells.Item(132, 9).Value2 = 1701.279  # I132: 1917.3243 -> 1701.279
$ws.Cells.Item(132, 10).Value2 = 4886.6665  # J132: 4967.273 -> 4886.6665
$ws.Cells.Item(132, 11).Value2 = 5103.837  # K132: 5751.9729 -> 5103.837
$ws.Cells.Item(132, 12).Value2 = 14659.9995  # L132: 14901.819 -> 14659.9995
$ws.Cells.Item(132, 13).Value2 = -2573.837  # M132: -3221.9729 -> -2573.837
$ws.Cells.Item(132, 14).Value2 = -19719.9995  # N132: -19961.819 -> -19719.9995
# Row 136
$ws.Cells.Item(136, 8).Value2 = 1593.7368  # H136: 66668280 -> 1593.7368
$ws.Cells.Item(136, 9).Value2 = 1493.7142  # I136: 76924510 -> 1493.7142
$ws.Cells.Item(136, 10).Value2 = 1873.8  # J136: 2749.5 -> 1873.8
$ws.Cells.Item(136, 11).Value2 = 4481.142599999999  # K136: 230773530 -> 4481.142599999999
$ws.Cells.Item(136, 12).Value2 = 5621.4  # L136: 8248.5 -> 5621.4
$ws.Cells.Item(136, 13).Value2 = -1931.142599999999  # M136: -230770980 -> -1931.142599999999
$ws.Cells.Item(136, 14).Value2 = -10721.4  # N136: -13348.5 -> -10721.4

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value2 = 2766.4614  # H20: 3151.111 -> 2766.4614
$ws.Cells.Item(20, 9).Value2 = 2833.889  # I20: 3316.6667 -> 2833.889
$ws.Cells.Item(20, 10).Value2 = 2614.75  # J20: 2820 -> 2614.75
$ws.Cells.Item(20, 11).Value2 = 2833.889  # K20: 3316.6667 -> 2833.889
$ws.Cells.Item(20, 12).Value2 = 2614.75  # L20: 2820 -> 2614.75
$ws.Cells.Item(20, 13).Value2 = -2586.889  # M20: -3069.6667 -> -2586.889
$ws.Cells.Item(20, 14).Value2 = -3108.75  # N20: -3314 -> -3108.75

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value2 = 52632676  # H16: 58824704 -> 52632676
$ws.Cells.Item(16, 9).Value2 = 71429690  # I16: 83334560 -> 71429690
$ws.Cells.Item(16, 11).Value2 = 71429690  # K16: 83334560 -> 71429690
$ws.Cells.Item(16, 13).Value2 = -71429403  # M16: -83334273 -> -71429403
# Row 31
$ws.Cells.Item(31, 8).Value2 = 1584.1207  # H31: 1625.2727 -> 1584.1207
$ws.Cells.Item(31, 9).Value2 = 1441.5918  # I31: 1475.5106 -> 1441.5918
$ws.Cells.Item(31, 10).Value2 = 2360.111  # J31: 2505.125 -> 2360.111
$ws.Cells.Item(31, 11).Value2 = 1441.5918  # K31: 1475.5106 -> 1441.5918
$ws.Cells.Item(31, 12).Value2 = 2360.111  # L31: 2505.125 -> 2360.111
$ws.Cells.Item(31, 13).Value2 = -1146.5918  # M31: -1180.5106 -> -1146.5918
$ws.Cells.Item(31, 14).Value2 = -2950.111  # N31: -3095.125 -> -2950.111
# Row 34
$ws.Cells.Item(34, 8).Value2 = 1584.1207  # H34: 1625.2727 -> 1584.1207
$ws.Cells.Item(34, 9).Value2 = 1441.5918  # I34: 1475.5106 -> 1441.5918
$ws.Cells.Item(34, 10).Value2 = 2360.111  # J34: 2505.125 -> 2360.111
$ws.Cells.Item(34, 11).Value2 = 1441.5918  # K34: 1475.5106 -> 1441.5918
$ws.Cells.Item(34, 12).Value2 = 2360.111  # L34: 2505.125 -> 2360.111
$ws.Cells.Item(34, 13).Value2 = -1239.5918  # M34: -1273.5106 -> -1239.5918
$ws.Cells.Item(34, 14).Value2 = -2764.111  # N34: -2909.125 -> -2764.111
# Row 58
$ws.Cells.Item(58, 8).Value2 = 4125.744  # H58: 4424.927 -> 4125.744
$ws.Cells.Item(58, 9).Value2 = 1150.2  # I58: 1219.4783 -> 1150.2
$ws.Cells.Item(58, 10).Value2 = 8258.444  # J58: 8520.777 -> 8258.444
$ws.Cells.Item(58, 11).Value2 = 1150.2  # K58: 1219.4783 -> 1150.2
$ws.Cells.Item(58, 12).Value2 = 8258.444  # L58: 8520.777 -> 8258.444
$ws.Cells.Item(58, 13).Value2 = -947.2  # M58: -1016.4783 -> -947.2
$ws.Cells.Item(58, 14).Value2 = -8664.444  # N58: -8926.777 -> -8664.444
# Row 99
$ws.Cells.Item(99, 8).Value2 = 1920  # H99: 1802.0714 -> 1920
$ws.Cells.Item(99, 9).Value2 = 1920  # I99: 1682.9 -> 1920
$ws.Cells.Item(99, 10).Value2 = 0  # J99: 2100 -> 0
$ws.Cells.Item(99, 11).Value2 = 1920  # K99: 1682.9 -> 1920
$ws.Cells.Item(99, 12).Value2 = 0  # L99: 2100 -> 0
$ws.Cells.Item(99, 13).Value2 = -422  # M99: -184.9000000000001 -> -422
$ws.Cells.Item(99, 14).ClearContents()  # N99: was -5096
# Row 113
$ws.Cells.Item(113, 8).Value2 = 52632676  # H113: 58824704 -> 52632676
$ws.Cells.Item(113, 9).Value2 = 71429690  # I113: 83334560 -> 71429690
$ws.Cells.Item(113, 11).Value2 = 71429690  # K113: 83334560 -> 71429690
$ws.Cells.Item(113, 13).Value2 = -71427520  # M113: -83332390 -> -71427520
# Row 122
$ws.Cells.Item(122, 8).Value2 = 1292  # H122: 920.93335 -> 1292
$ws.Cells.Item(122, 9).Value2 = 1156  # I122: 920.93335 -> 1156
$ws.Cells.Item(122, 10).Value2 = 1700  # J122: 0 -> 1700
$ws.Cells.Item(122, 11).Value2 = 3468  # K122: 2762.80005 -> 3468
$ws.Cells.Item(122, 12).Value2 = 5100  # L122: 0 -> 5100
$ws.Cells.Item(122, 13).Value2 = -1018  # M122: -312.8000499999998 -> -1018
$ws.Cells.Item(122, 14).Value2 = -10000  # N122: None -> -10000
# Row 126
$ws.Cells.Item(126, 8).Value2 = 1920  # H126: 1802.0714 -> 1920
$ws.Cells.Item(126, 9).Value2 = 1920  # I126: 1682.9 -> 1920
$ws.Cells.Item(126, 10).Value2 = 0  # J126: 2100 -> 0
$ws.Cells.Item(126, 11).Value2 = 5760  # K126: 5048.700000000001 -> 5760
$ws.Cells.Item(126, 12).Value2 = 0  # L126: 6300 -> 0
$ws.Cells.Item(126, 13).Value2 = -3290  # M126: -2578.700000000001 -> -3290
$ws.Cells.Item(126, 14).ClearContents()  # N126: was -11240
# Row 134
$ws.Cells.Item(134, 8).Value2 = 1619.6428  # H134: 16668205 -> 1619.6428
$ws.Cells.Item(134, 9).Value2 = 1615.3182  # I134: 1579.8695 -> 1615.3182
$ws.Cells.Item(134, 10).Value2 = 1635.5  # J134: 71429976 -> 1635.5
$ws.Cells.Item(134, 11).Value2 = 4845.9546  # K134: 4739.6085 -> 4845.9546
$ws.Cells.Item(134, 12).Value2 = 4906.5  # L134: 214289928 -> 4906.5
$ws.Cells.Item(134, 13).Value2 = -2310.9546  # M134: -2204.6085 -> -2310.9546
$ws.Cells.Item(134, 14).Value2 = -9976.5  # N134: -214294998 -> -9976.5
# Row 136
$ws.Cells.Item(136, 8).Value2 = 4125.744  # H136: 4424.927 -> 4125.744
$ws.Cells.Item(136, 9).Value2 = 1150.2  # I136: 1219.4783 -> 1150.2
$ws.Cells.Item(136, 10).Value2 = 8258.444  # J136: 8520.777 -> 8258.444
$ws.Cells.Item(136, 11).Value2 = 3450.6  # K136: 3658.4349 -> 3450.6
$ws.Cells.Item(136, 12).Value2 = 24775.332  # L136: 25562.331 -> 24775.332
$ws.Cells.Item(136, 13).Value2 = -900.6000000000004  # M136: -1108.4349 -> -900.6000000000004
$ws.Cells.Item(136, 14).Value2 = -29875.332  # N136: -30662.331 -> -29875.332

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value2 = 741.8461  # H68: 747.2308 -> 741.8461
$ws.Cells.Item(68, 10).Value2 = 589.3333  # J68: 597.1111 -> 589.3333
$ws.Cells.Item(68, 12).Value2 = 1767.9999  # L68: 1791.3333 -> 1767.9999
$ws.Cells.Item(68, 14).Value2 = -3389.9999  # N68: -3413.3333 -> -3389.9999
# Row 71
$ws.Cells.Item(71, 8).Value2 = 741.8461  # H71: 747.2308 -> 741.8461
$ws.Cells.Item(71, 10).Value2 = 589.3333  # J71: 597.1111 -> 589.3333
$ws.Cells.Item(71, 12).Value2 = 5303.9997  # L71: 5373.9999 -> 5303.9997
$ws.Cells.Item(71, 14).Value2 = -13415.9997  # N71: -13485.9999 -> -13415.9997
# Row 80
$ws.Cells.Item(80, 8).Value2 = 1500  # H80: 0 -> 1500
$ws.Cells.Item(80, 10).Value2 = 1500  # J80: 0 -> 1500
$ws.Cells.Item(80, 12).Value2 = 4500  # L80: 0 -> 4500
$ws.Cells.Item(80, 14).Value2 = -6372  # N80: None -> -6372
# Row 83
$ws.Cells.Item(83, 8).Value2 = 1500  # H83: 0 -> 1500
$ws.Cells.Item(83, 10).Value2 = 1500  # J83: 0 -> 1500
$ws.Cells.Item(83, 12).Value2 = 13500  # L83: 0 -> 13500
$ws.Cells.Item(83, 14).Value2 = -22860  # N83: None -> -22860
# Row 100
$ws.Cells.Item(100, 8).Value2 = 13268.777  # H100: 3236.3333 -> 13268.777
$ws.Cells.Item(100, 10).Value2 = 13268.777  # J100: 3236.3333 -> 13268.777
$ws.Cells.Item(100, 12).Value2 = 39806.331  # L100: 9708.999899999999 -> 39806.331
$ws.Cells.Item(100, 14).Value2 = -41428.331  # N100: -11330.9999 -> -41428.331
# Row 131
$ws.Cells.Item(131, 8).Value2 = 18547158  # H131: 21309402 -> 18547158
$ws.Cells.Item(131, 9).Value2 = 58824028  # I131: 62500530 -> 58824028
$ws.Cells.Item(131, 10).Value2 = 41569.676  # J131: 49464.453 -> 41569.676
$ws.Cells.Item(131, 11).Value2 = 176472084  # K131: 187501590 -> 176472084
$ws.Cells.Item(131, 12).Value2 = 124709.028  # L131: 148393.359 -> 124709.028
$ws.Cells.Item(131, 13).Value2 = -176467044  # M131: -187496550 -> -176467044
$ws.Cells.Item(131, 14).Value2 = -134789.028  # N131: -158473.359 -> -134789.028
# Row 132
$ws.Cells.Item(132, 8).Value2 = 1047.6666  # H132: 1206.0714 -> 1047.6666
$ws.Cells.Item(132, 9).Value2 = 796.38464  # I132: 931.1111 -> 796.38464
$ws.Cells.Item(132, 11).Value2 = 7167.46176  # K132: 8379.999899999999 -> 7167.46176
$ws.Cells.Item(132, 13).Value2 = -4637.46176  # M132: -5849.999899999999 -> -4637.46176
# Row 137
$ws.Cells.Item(137, 8).Value2 = 34098024  # H137: 32615790 -> 34098024
$ws.Cells.Item(137, 9).Value2 = 53573388  # I137: 57694270 -> 53573388
$ws.Cells.Item(137, 10).Value2 = 16139.375  # J137: 13765.5 -> 16139.375
$ws.Cells.Item(137, 11).Value2 = 160720164  # K137: 173082810 -> 160720164
$ws.Cells.Item(137, 12).Value2 = 48418.125  # L137: 41296.5 -> 48418.125
$ws.Cells.Item(137, 13).Value2 = -160715064  # M137: -173077710 -> -160715064
$ws.Cells.Item(137, 14).Value2 = -58618.125  # N137: -51496.5 -> -58618.125

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value2 = 319.8889  # H2: 310.60715 -> 319.8889
$ws.Cells.Item(2, 9).Value2 = 239.73334  # I2: 228.5 -> 239.73334
$ws.Cells.Item(2, 11).Value2 = 239.73334  # K2: 228.5 -> 239.73334
$ws.Cells.Item(2, 13).Value2 = -126.73334  # M2: -115.5 -> -126.73334
# Row 26
$ws.Cells.Item(26, 8).Value2 = 0  # H26: 25000 -> 0
$ws.Cells.Item(26, 10).Value2 = 0  # J26: 25000 -> 0
$ws.Cells.Item(26, 12).Value2 = 0  # L26: 25000 -> 0
$ws.Cells.Item(26, 14).ClearContents()  # N26: was -25560
# Row 50
$ws.Cells.Item(50, 8).Value2 = 0  # H50: 25000 -> 0
$ws.Cells.Item(50, 10).Value2 = 0  # J50: 25000 -> 0
$ws.Cells.Item(50, 12).Value2 = 0  # L50: 25000 -> 0
$ws.Cells.Item(50, 14).ClearContents()  # N50: was -25996
# Row 70
$ws.Cells.Item(70, 8).Value2 = 40913080  # H70: 40913036 -> 40913080
$ws.Cells.Item(70, 9).Value2 = 31254300  # I70: 27782044 -> 31254300
$ws.Cells.Item(70, 10).Value2 = 66669830  # J70: 100002500 -> 66669830
$ws.Cells.Item(70, 11).Value2 = 31254300  # K70: 27782044 -> 31254300
$ws.Cells.Item(70, 12).Value2 = 66669830  # L70: 100002500 -> 66669830
$ws.Cells.Item(70, 13).Value2 = -31254030  # M70: -27781774 -> -31254030
$ws.Cells.Item(70, 14).Value2 = -66670370  # N70: -100003040 -> -66670370
# Row 73
$ws.Cells.Item(73, 8).Value2 = 40913080  # H73: 40913036 -> 40913080
$ws.Cells.Item(73, 9).Value2 = 31254300  # I73: 27782044 -> 31254300
$ws.Cells.Item(73, 10).Value2 = 66669830  # J73: 100002500 -> 66669830
$ws.Cells.Item(73, 11).Value2 = 31254300  # K73: 27782044 -> 31254300
$ws.Cells.Item(73, 12).Value2 = 66669830  # L73: 100002500 -> 66669830
$ws.Cells.Item(73, 13).Value2 = -31253364  # M73: -27781108 -> -31253364
$ws.Cells.Item(73, 14).Value2 = -66671702  # N73: -100004372 -> -66671702
# Row 97
$ws.Cells.Item(97, 8).Value2 = 636.0454999999999  # H97: 662.5714 -> 636.0454999999999
$ws.Cells.Item(97, 9).Value2 = 538.5  # I97: 565.5294 -> 538.5
$ws.Cells.Item(97, 11).Value2 = 538.5  # K97: 565.5294 -> 538.5
$ws.Cells.Item(97, 13).Value2 = -42.5  # M97: -69.52940000000001 -> -42.5
# Row 132
$ws.Cells.Item(132, 8).Value2 = 2836.7317  # H132: 3071.9487 -> 2836.7317
$ws.Cells.Item(132, 9).Value2 = 2565.4375  # I132: 2892.8965 -> 2565.4375
$ws.Cells.Item(132, 10).Value2 = 3801.3333  # J132: 3591.2 -> 3801.3333
$ws.Cells.Item(132, 11).Value2 = 7696.3125  # K132: 8678.6895 -> 7696.3125
$ws.Cells.Item(132, 12).Value2 = 11403.9999  # L132: 10773.6 -> 11403.9999
$ws.Cells.Item(132, 13).Value2 = -5166.3125  # M132: -6148.6895 -> -5166.3125
$ws.Cells.Item(132, 14).Value2 = -16463.9999  # N132: -15833.6 -> -16463.9999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Cells.Item(132, 8).Value2 = 2539.6428  # H132: 2604.423 -> 2539.6428
$ws.Cells.Item(132, 9).Value2 = 2607.6155  # I132: 2600.2307 -> 2607.6155
$ws.Cells.Item(132, 10).Value2 = 2480.7334  # J132: 2608.6155 -> 2480.7334
$ws.Cells.Item(132, 11).Value2 = 7822.8465  # K132: 7800.6921 -> 7822.8465
$ws.Cells.Item(132, 12).Value2 = 7442.2002  # L132: 7825.8465 -> 7442.2002
$ws.Cells.Item(132, 13).Value2 = -5292.8465  # M132: -5270.6921 -> -5292.8465
$ws.Cells.Item(132, 14).Value2 = -12502.2002  # N132: -12885.8465 -> -12502.2002

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Cells.Item(51, 8).Value2 = 0  # H51: 11075 -> 0
$ws.Cells.Item(51, 10).Value2 = 0  # J51: 11075 -> 0
$ws.Cells.Item(51, 12).Value2 = 0  # L51: 11075 -> 0
$ws.Cells.Item(51, 14).ClearContents()  # N51: was -12095
# Row 132
$ws.Cells.Item(132, 8).Value2 = 2205.6  # H132: 2343.9688 -> 2205.6
$ws.Cells.Item(132, 9).Value2 = 1757.25  # I132: 1880.56 -> 1757.25
$ws.Cells.Item(132, 11).Value2 = 5271.75  # K132: 5641.68 -> 5271.75
$ws.Cells.Item(132, 13).Value2 = -2741.75  # M132: -3111.68 -> -2741.75
# Row 136
$ws.Cells.Item(136, 8).Value2 = 1941.3572  # H136: 1866.2 -> 1941.3572
$ws.Cells.Item(136, 9).Value2 = 1431  # I136: 1369.3 -> 1431
$ws.Cells.Item(136, 11).Value2 = 4293  # K136: 4107.9 -> 4293
$ws.Cells.Item(136, 13).Value2 = -1743  # M136: -1557.9 -> -1743
